# Add a new row "A21 / Pretraga arhive" right after the "A2 / Pregled arhive"
# row in the (single) table of the document.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Locate the row that contains "Pregled arhive" in its second column so the
# new row can be inserted directly below it (mirrors the XML diff, where the
# new row sits between the A2 and A3 rows).
$targetRowIndex = -1
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $cellText = $table.Cell($i, 2).Range.Text
    if ($cellText -like "Pregled arhive*") {
        $targetRowIndex = $i
        break
    }
}

if ($targetRowIndex -eq -1) {
    throw "Could not find the 'Pregled arhive' row"
}

$insertBeforeIndex = $targetRowIndex + 1
$insertBeforeRow = $table.Rows.Item($insertBeforeIndex)
$newRow = $table.Rows.Add($insertBeforeRow)

$newRowIndex = $targetRowIndex + 1
$table.Cell($newRowIndex, 1).Range.Text = "A21"
$table.Cell($newRowIndex, 2).Range.Text = "Pretraga arhive"
